# Update the "last updated" timestamp shown in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 16:52"

# Row 4 - Estados Unidos: refreshed case counts
$ws.Cells.Item(4, 2).Value = 928364
$ws.Cells.Item(4, 3).Value = 3132
$ws.Cells.Item(4, 4).Value = 110490
$ws.Cells.Item(4, 5).Value = 765518
$ws.Cells.Item(4, 7).Value = 163
$ws.Cells.Item(4, 8).Value = 52356

# Row 8 - refreshed case counts
$ws.Cells.Item(8, 2).Value = 155407
$ws.Cells.Item(8, 3).Value = 408
$ws.Cells.Item(8, 5).Value = 39805
$ws.Cells.Item(8, 7).Value = 42
$ws.Cells.Item(8, 8).Value = 5802

# Row 9 - refreshed case counts
$ws.Cells.Item(9, 2).Value = 148377
$ws.Cells.Item(9, 3).Value = 4913
$ws.Cells.Item(9, 5).Value = 127714
$ws.Cells.Item(9, 7).Value = 813
$ws.Cells.Item(9, 8).Value = 20319

# Row 57 - Argentina: refreshed case counts
$ws.Cells.Item(57, 4).Value = 1030
$ws.Cells.Item(57, 5).Value = 2398
$ws.Cells.Item(57, 6).Value = 144
$ws.Cells.Item(57, 7).Value = 3
$ws.Cells.Item(57, 8).Value = 179

# Rows 58/59 - Moldavia overtakes Argelia in ranking, so the two rows swap
# (Moldavia moves up to row 58 with refreshed counts, Argelia moves down to
# row 59 keeping its previous counts)
$ws.Cells.Item(58, 1).Value = "Moldavia"
$ws.Cells.Item(58, 2).Value = 3304
$ws.Cells.Item(58, 3).Value = 194
$ws.Cells.Item(58, 4).Value = 825
$ws.Cells.Item(58, 5).Value = 2385
$ws.Cells.Item(58, 6).Value = 212
$ws.Cells.Item(58, 7).Value = 10
$ws.Cells.Item(58, 8).Value = 94

$ws.Cells.Item(59, 1).Value = "Argelia"
$ws.Cells.Item(59, 2).Value = 3127
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 1408
$ws.Cells.Item(59, 5).Value = 1304
$ws.Cells.Item(59, 6).Value = 40
$ws.Cells.Item(59, 7).Value = 0
$ws.Cells.Item(59, 8).Value = 415

# Row 84 - refreshed case counts
$ws.Cells.Item(84, 2).Value = 1247
$ws.Cells.Item(84, 3).Value = 59
$ws.Cells.Item(84, 5).Value = 995
$ws.Cells.Item(84, 7).Value = 1
$ws.Cells.Item(84, 8).Value = 55

# Row 115 - refreshed case counts
$ws.Cells.Item(115, 5).Value = 364
$ws.Cells.Item(115, 7).Value = 2
$ws.Cells.Item(115, 8).Value = 18
